$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1124
$ws.Range("J32").Value = 1166
$ws.Range("L32").Value = 1166
$ws.Range("N32").Value = -1818
$ws.Range("H55").Value = 323.57144
$ws.Range("I55").Value = 301.81818
$ws.Range("J55").Value = 403.33334
$ws.Range("K55").Value = 301.81818
$ws.Range("L55").Value = 403.33334
$ws.Range("M55").Value = -87.81817999999998
$ws.Range("N55").Value = -831.33334
$ws.Range("H92").Value = 693.75
$ws.Range("I92").Value = 705.1177
$ws.Range("K92").Value = 705.1177
$ws.Range("M92").Value = 542.8823
$ws.Range("H100").Value = 3244.6667
$ws.Range("I100").Value = 1708.75
$ws.Range("K100").Value = 1708.75
$ws.Range("M100").Value = -1167.75
$ws.Range("H135").Value = 1711.6923
$ws.Range("I135").Value = 1333
$ws.Range("K135").Value = 11997
$ws.Range("M135").Value = -9462
$ws.Range("H137").Value = 3705277.2
$ws.Range("I137").Value = 1184.4286
$ws.Range("K137").Value = 3553.2858
$ws.Range("M137").Value = -1003.2858

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22728816
$ws.Range("I32").Value = 25001264
$ws.Range("K32").Value = 25001264
$ws.Range("M32").Value = -25000977
$ws.Range("H44").Value = 68993
$ws.Range("J44").Value = 68986
$ws.Range("L44").Value = 68986
$ws.Range("N44").Value = -69962
$ws.Range("H45").Value = 4872.148
$ws.Range("I45").Value = 6425
$ws.Range("J45").Value = 4218.316
$ws.Range("K45").Value = 6425
$ws.Range("L45").Value = 4218.316
$ws.Range("M45").Value = -6048
$ws.Range("N45").Value = -4972.316
$ws.Range("H55").Value = 40021
$ws.Range("J55").Value = 69994
$ws.Range("L55").Value = 69994
$ws.Range("N55").Value = -70624
$ws.Range("H74").Value = 2986.111
$ws.Range("I74").Value = 2986.111
$ws.Range("K74").Value = 2986.111
$ws.Range("M74").Value = -2112.111
$ws.Range("H77").Value = 2986.111
$ws.Range("I77").Value = 2986.111
$ws.Range("K77").Value = 14930.555
$ws.Range("M77").Value = -10562.555
$ws.Range("H132").Value = 4928.8125
$ws.Range("I132").Value = 5258.846
$ws.Range("K132").Value = 15776.538
$ws.Range("M132").Value = -13246.538

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29916.777
$ws.Range("I20").Value = 47200.727
$ws.Range("J20").Value = 2756.2856
$ws.Range("K20").Value = 47200.727
$ws.Range("L20").Value = 2756.2856
$ws.Range("M20").Value = -46953.727
$ws.Range("N20").Value = -3250.2856
$ws.Range("H86").Value = 6177.8
$ws.Range("J86").Value = 6823
$ws.Range("L86").Value = 6823
$ws.Range("N86").Value = -9069
$ws.Range("H89").Value = 6177.8
$ws.Range("J89").Value = 6823
$ws.Range("L89").Value = 34115
$ws.Range("N89").Value = -45347
$ws.Range("H94").Value = 697.64703
$ws.Range("I94").Value = 752.9286
$ws.Range("K94").Value = 752.9286
$ws.Range("M94").Value = -301.9286
$ws.Range("H107").Value = 2639.7273
$ws.Range("I107").Value = 2837.7778
$ws.Range("J107").Value = 1748.5
$ws.Range("K107").Value = 2837.7778
$ws.Range("L107").Value = 1748.5
$ws.Range("M107").Value = -917.7777999999998
$ws.Range("N107").Value = -5588.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1305
$ws.Range("I16").Value = 1305
$ws.Range("K16").Value = 1305
$ws.Range("M16").Value = -1018
$ws.Range("H31").Value = 5721.84
$ws.Range("J31").Value = 5994.8184
$ws.Range("L31").Value = 5994.8184
$ws.Range("N31").Value = -6584.8184
$ws.Range("H34").Value = 5721.84
$ws.Range("J34").Value = 5994.8184
$ws.Range("L34").Value = 5994.8184
$ws.Range("N34").Value = -6398.8184
$ws.Range("H86").Value = 19003.592
$ws.Range("I86").Value = 19491.066
$ws.Range("K86").Value = 19491.066
$ws.Range("M86").Value = -18368.066
$ws.Range("H89").Value = 19003.592
$ws.Range("I89").Value = 19491.066
$ws.Range("K89").Value = 97455.32999999999
$ws.Range("M89").Value = -91839.32999999999
$ws.Range("H99").Value = 3999.5
$ws.Range("I99").Value = 3999.5
$ws.Range("K99").Value = 3999.5
$ws.Range("M99").Value = -2501.5
$ws.Range("H107").Value = 1278.303
$ws.Range("I107").Value = 798.625
$ws.Range("J107").Value = 1431.8
$ws.Range("K107").Value = 798.625
$ws.Range("L107").Value = 1431.8
$ws.Range("M107").Value = 1121.375
$ws.Range("N107").Value = -5271.8
$ws.Range("H113").Value = 1305
$ws.Range("I113").Value = 1305
$ws.Range("K113").Value = 1305
$ws.Range("M113").Value = 865
$ws.Range("H122").Value = 4834.6875
$ws.Range("I122").Value = 4041.111
$ws.Range("J122").Value = 5855
$ws.Range("K122").Value = 12123.333
$ws.Range("L122").Value = 17565
$ws.Range("M122").Value = -9673.332999999999
$ws.Range("N122").Value = -22465
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 3999.5
$ws.Range("I126").Value = 3999.5
$ws.Range("K126").Value = 11998.5
$ws.Range("M126").Value = -9528.5
$ws.Range("H132").Value = 3022.818
$ws.Range("I132").Value = 2581.889
$ws.Range("K132").Value = 7745.667
$ws.Range("M132").Value = -5215.667
$ws.Range("H134").Value = 2501
$ws.Range("I134").Value = 1996.6666
$ws.Range("K134").Value = 5989.9998
$ws.Range("M134").Value = -3454.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 466.66666
$ws.Range("J50").Value = 466.66666
$ws.Range("L50").Value = 1399.99998
$ws.Range("N50").Value = -2361.99998
$ws.Range("H53").Value = 466.66666
$ws.Range("J53").Value = 466.66666
$ws.Range("L53").Value = 1399.99998
$ws.Range("N53").Value = -2361.99998
$ws.Range("H121").Value = 938.5
$ws.Range("J121").Value = 2198.6667
$ws.Range("L121").Value = 6596.000100000001
$ws.Range("N121").Value = -9216.000100000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 21371.771
$ws.Range("I113").Value = 8645
$ws.Range("J113").Value = 72278.86
$ws.Range("K113").Value = 8645
$ws.Range("L113").Value = 72278.86
$ws.Range("M113").Value = -6475
$ws.Range("N113").Value = -76618.86
$ws.Range("H122").Value = 1952.3
$ws.Range("I122").Value = 2047.5714
$ws.Range("J122").Value = 1730
$ws.Range("K122").Value = 6142.7142
$ws.Range("L122").Value = 5190
$ws.Range("M122").Value = -3692.7142
$ws.Range("N122").Value = -10090
$ws.Range("H126").Value = 2444.2273
$ws.Range("I126").Value = 2058.4
$ws.Range("J126").Value = 2765.75
$ws.Range("K126").Value = 6175.200000000001
$ws.Range("L126").Value = 8297.25
$ws.Range("M126").Value = -3705.200000000001
$ws.Range("N126").Value = -13237.25
$ws.Range("H132").Value = 2012
$ws.Range("I132").Value = 2012
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6036
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3506
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3882.3635
$ws.Range("I7").Value = 3819.4285
$ws.Range("J7").Value = 3992.5
$ws.Range("K7").Value = 3819.4285
$ws.Range("L7").Value = 3992.5
$ws.Range("M7").Value = -3707.4285
$ws.Range("N7").Value = -4216.5
$ws.Range("H16").Value = 1371.9565
$ws.Range("I16").Value = 1431.3334
$ws.Range("J16").Value = 748.5
$ws.Range("K16").Value = 1431.3334
$ws.Range("L16").Value = 748.5
$ws.Range("M16").Value = -1261.3334
$ws.Range("N16").Value = -1088.5
$ws.Range("H40").Value = 2652.111
$ws.Range("I40").Value = 4882
$ws.Range("J40").Value = 2015
$ws.Range("K40").Value = 4882
$ws.Range("L40").Value = 2015
$ws.Range("M40").Value = -4746
$ws.Range("N40").Value = -2287
$ws.Range("H46").Value = 3543.6128
$ws.Range("I46").Value = 2613
$ws.Range("K46").Value = 2613
$ws.Range("M46").Value = -2425
$ws.Range("H126").Value = 3882.3635
$ws.Range("I126").Value = 3819.4285
$ws.Range("J126").Value = 3992.5
$ws.Range("K126").Value = 11458.2855
$ws.Range("L126").Value = 11977.5
$ws.Range("M126").Value = -8988.2855
$ws.Range("N126").Value = -16917.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 28518
$ws.Range("I53").Value = 30036
$ws.Range("J53").Value = 27000
$ws.Range("K53").Value = 30036
$ws.Range("L53").Value = 27000
$ws.Range("M53").Value = -29429
$ws.Range("N53").Value = -28214
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H122").Value = 4563.543
$ws.Range("I122").Value = 4882.25
$ws.Range("J122").Value = 3868.182
$ws.Range("K122").Value = 14646.75
$ws.Range("L122").Value = 11604.546
$ws.Range("M122").Value = -12196.75
$ws.Range("N122").Value = -16504.546
$ws.Range("H126").Value = 3299.6667
$ws.Range("I126").Value = 2683.6667
$ws.Range("J126").Value = 4531.6665
$ws.Range("K126").Value = 8051.000100000001
$ws.Range("L126").Value = 13594.9995
$ws.Range("M126").Value = -5581.000100000001
$ws.Range("N126").Value = -18534.9995
$ws.Range("H136").Value = 2245.64
$ws.Range("I136").Value = 1769.7273
$ws.Range("K136").Value = 5309.1819
$ws.Range("M136").Value = -2759.1819

Write-Output "applied 243 cell updates across 8 sheets"